$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Add the new log entry (Post 41) in row 51 ---
# Set the string-valued cells in the same order the source workbook was
# authored in (Hashnode link, then Title, then Dev.to link) so new
# shared-string table entries land in the expected order.
$ws.Range("E51").Value = "https://programmingport.hashnode.dev/file-test-operator-or-shell-scripting"
$ws.Range("C51").Value = "File Test Operator | Shell Scripting"
$ws.Range("F51").Value = "https://dev.to/rahulmishra05/file-test-operator-shell-scripting-2hha"
$ws.Range("B51").Value = 41
$ws.Range("D51").Value = 44169

# Match the formatting used by the rest of the table body (row 50):
# plain numbers/text for S.No & Title, a date format for Date of Post,
# and the Hyperlink cell style for the two link columns.
$ws.Range("D51").NumberFormat = $ws.Range("D50").NumberFormat
$ws.Range("E51").Style = "Hyperlink"
$ws.Range("F51").Style = "Hyperlink"

# --- Grow Table2 so the new row becomes part of the table ---
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("B10:F51"))

# --- Update the view so the newly-added row/cell is what's in focus ---
$excel.Goto($ws.Range("F51"), $true)
